$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$values = @{
    2  = @(2540.05, 2519.8)
    3  = @(376.4,   382.55)
    4  = @(1538.1,  1550.9)
    5  = @(7492.2,  7494.3)
    6  = @(249.35,  239.55)
    7  = @(210.8,   211.9)
    8  = @(46288.05,46113.65)
    9  = @(532.4,   537.35)
    10 = @(3427.2,  3381.4)
    11 = @(148.25,  147.5)
    12 = @(1304.6,  1283.6)
    13 = @(1450,    1458.25)
    14 = @(711.45,  711.05)
    15 = @(458.75,  461.75)
    16 = @(1601.1,  1566)
    17 = @(300.95,  300.15)
    18 = @(20249.6, 20172)
    19 = @(598.8,   597.3)
    20 = @(596.2,   605.45)
    21 = @(634.25,  624.6)
    22 = @(263.95,  265.1)
    23 = @(131.95,  131.7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
